$d = $word.ActiveDocument

$map = @(
    @("987×3=", "122×2="),
    @("121×7=", "868×9="),
    @("849×5=", "327×8="),
    @("513×8=", "200×2="),
    @("564×4=", "116×7="),
    @("186×8=", "685×9="),
    @("439×8=", "891×3="),
    @("392×2=", "603×7="),
    @("420×5=", "964×5="),
    @("295×6=", "322×8="),
    @("982×2=", "681×5="),
    @("449×9=", "587×8="),
    @("298×3=", "325×8="),
    @("491×4=", "602×7="),
    @("252×6=", "481×2="),
    @("542×5=", "128×5="),
    @("121×6=", "125×8="),
    @("800×9=", "771×5="),
    @("428×6=", "532×7="),
    @("927×6=", "835×3="),
    @("267×8=", "760×5="),
    @("193×6=", "142×8="),
    @("393×8=", "112×8="),
    @("214×8=", "770×3="),
    @("947×9=", "236×3=")
)

foreach ($pair in $map) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
